# Apply "Automatic update of files" changes to the Avverkningsanmälningar sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") on rows 2-20: date serial 45208 -> 45212 (2023-10-09 -> 2023-10-13)
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}

# 2) Rows 2-4 (cases A 30840-2023, A 30841-2023, A 30839-2023): the hyperlink formulas in
#    columns S-Y get extra descriptive suffixes appended to the linked filenames.
#    Column -> (folder name in URL, file suffix to insert before the extension, file extension)
$linkCols = @(
    @{ Col = "S"; Folder = "artfynd";        Suffix = " artfynd";               Ext = "xlsx" },
    @{ Col = "T"; Folder = "kartor";         Suffix = " karta";                 Ext = "png"  },
    @{ Col = "U"; Folder = "knärot";         Suffix = " karta knärot";          Ext = "png"  },
    @{ Col = "V"; Folder = "klagomål";       Suffix = " fsc-klagomål";          Ext = "docx" },
    @{ Col = "W"; Folder = "klagomålsmail";  Suffix = " fsc-klagomål mail";     Ext = "docx" },
    @{ Col = "X"; Folder = "tillsyn";        Suffix = " tillsynsbegäran";       Ext = "docx" },
    @{ Col = "Y"; Folder = "ti,llsynsmail";  Suffix = " tillsynsbegäran mail"; Ext = "docx" }
)

$caseNames = @{
    2 = "A 30840-2023"
    3 = "A 30841-2023"
    4 = "A 30839-2023"
}

foreach ($row in 2..4) {
    $name = $caseNames[$row]
    foreach ($link in $linkCols) {
        $cell = $ws.Range("$($link.Col)$row")
        # Only touch cells that already contain a HYPERLINK formula (e.g. U4 is blank).
        $existing = $cell.Formula
        if ($existing -like "*HYPERLINK*") {
            $url = "https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/$($link.Folder)/$name$($link.Suffix).$($link.Ext)"
            $cell.Formula = '=HYPERLINK("' + $url + '", "' + $name + '")'
        }
    }
}
